$d = $word.ActiveDocument
$para = $d.Paragraphs(1)
$range = $para.Range
$range.Collapse(1)
$range.InsertBefore("**Rework the database so that every artist/song has to be searched from Spotify. Then, upon object creation, store the Spotify object ID to facilitate lookup when someone logs in with Spotify to search for live performance of songs/artists that they like.**`r")
